$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.556.92"
$ws.Range("E2").Value = "'  +4.01%  "
$ws.Range("D3").Value = "'3.486.08"
$ws.Range("E3").Value = "'  +2.57%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'590.26"
$ws.Range("E5").Value = "'  +3.32%  "
$ws.Range("D6").Value = "'169.15"
$ws.Range("E6").Value = "'  +3.91%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("D8").Value = "'3.485.06"
$ws.Range("E8").Value = "'  +2.49%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "'  +7.85%  "
$ws.Range("D10").Value = "'7.32"
$ws.Range("E10").Value = "'  +0.40%  "
$ws.Range("E11").Value = "'  +6.28%  "
$ws.Range("D12").Value = "'0.438"
$ws.Range("E12").Value = "'  +4.08%  "
$ws.Range("D13").Value = "'4.088.83"
$ws.Range("E13").Value = "'  +2.58%  "
$ws.Range("E14").Value = "'  -0.07%  "
$ws.Range("D15").Value = "'28.11"
$ws.Range("E15").Value = "'  +4.43%  "
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("E16").Value = "'  +3.53%  "
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'66.591.83"
$ws.Range("E17").Value = "'  +3.99%  "
$ws.Range("D18").Value = "'3.495.99"
$ws.Range("E18").Value = "'  +2.81%  "
$ws.Range("E19").Value = "'  +3.18%  "
$ws.Range("D20").Value = "'14.05"
$ws.Range("E20").Value = "'  +3.41%  "
$ws.Range("D21").Value = "'391.01"
$ws.Range("E21").Value = "'  +3.87%  "
$ws.Range("D22").Value = "'7.89"
$ws.Range("E22").Value = "'  +1.57%  "
$ws.Range("D23").Value = "'73.02"
$ws.Range("E23").Value = "'  +3.90%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "'  -0.19%  "
$ws.Range("E25").Value = "'  +4.61%  "
$ws.Range("E26").Value = "'  +6.02%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("E27").Value = "'  +7.33%  "
$ws.Range("E28").Value = "'  +1.80%  "
$ws.Range("E29").Value = "'  +0.08%  "
$ws.Range("E30").Value = "'  +3.28%  "
$ws.Range("E31").Value = "'  +5.23%  "
$ws.Range("E32").Value = "'  +2.92%  "
$ws.Range("D33").Value = "'23.53"
$ws.Range("E33").Value = "'  +3.36%  "
$ws.Range("D34").Value = "'7.42"
$ws.Range("E34").Value = "'  +5.60%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "'  +0.09%  "
$ws.Range("D36").Value = "'1.62"
$ws.Range("E36").Value = "'  +9.26%  "
$ws.Range("D37").Value = "'161.81"
$ws.Range("E37").Value = "'  +1.26%  "
$ws.Range("D38").Value = "'0.888"
$ws.Range("E38").Value = "'  +3.23%  "
$ws.Range("E39").Value = "'  +5.96%  "
$ws.Range("B40").Value = "'Filecoin"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.65"
$ws.Range("E40").Value = "'  +6.46%  "
$ws.Range("B41").Value = "'RenderToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'6.72"
$ws.Range("E41").Value = "'  +4.75%  "
$ws.Range("B42").Value = "'Hedera"
$ws.Range("C42").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0741"
$ws.Range("E42").Value = "'  +2.99%  "
$ws.Range("D43").Value = "'26.44"
$ws.Range("E43").Value = "'  +2.78%  "
$ws.Range("D44").Value = "'26.70"
$ws.Range("E44").Value = "'  +1.81%  "
$ws.Range("D45").Value = "'43.08"
$ws.Range("E45").Value = "'  +0.60%  "
$ws.Range("D46").Value = "'2.766.25"
$ws.Range("E46").Value = "'  +0.90%  "
$ws.Range("E47").Value = "'  +1.99%  "
$ws.Range("E48").Value = "'  +2.87%  "
$ws.Range("D49").Value = "'346.12"
$ws.Range("E49").Value = "'  +5.73%  "
$ws.Range("E50").Value = "'  +4.81%  "
$ws.Range("D51").Value = "'0.883"
$ws.Range("E51").Value = "'  +8.71%  "